# katalog.xlsx - "Draht_Matten" sheet rework:
#  - Montageart option (row 5) changes from a "P_Fund" factor (Betonieren:9, Konsole:1)
#    to a plain on/off flag "Ist_Beton" (Betonieren:1, Konsole:0)
#  - a new row is inserted for the "Preis pro Sack Beton (€)" / P_Sack input
#  - the final price formula (row 9) is rewritten to use the new variables
#  - column B gets a fixed width so the longer labels are readable

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Draht_Matten")

# Insert a new row above the old "Konsolen" row (old row 6) to make room for
# the new "Preis pro Sack Beton (€)" input; everything below shifts down by one.
$ws.Rows.Item(6).Insert()

# Row 6 (new): price per bag of concrete.
$ws.Range("A6").Value = "Zahl"
$ws.Range("B6").Value = "Preis pro Sack Beton (€)"
$ws.Range("C6").Value = "P_Sack"

# Row 5 ("Montageart"): swap the P_Fund multiplier for a simple Ist_Beton flag.
$ws.Range("C5").Value = "Ist_Beton"
$ws.Range("D5").Value = "Betonieren:1, Konsole:0"

# Row 7 (was row 6, "Konsolen" - content unchanged, just shifted down already).
# Row 8 (was row 7, "Montage (€/m)" - content unchanged, just shifted down already).

# Row 9 (was row 8): updated final formula using the new variables.
$ws.Range("E9").Value = "(L * P_Matte_Lfm * F_Faktor) + ((math.ceil(L/2.5)+1) * ( (P_Saeule * F_Faktor) + (Ist_Beton * 2 * P_Sack) + ((1-Ist_Beton) * P_Konsole * F_Faktor) )) + (Ecken * 30) + (L * P_Arbeit)"

# Widen column B to fit the new, longer labels.
$ws.Columns.Item(2).ColumnWidth = 20.83

# Match the author's final cursor position when they saved.
$ws.Range("E11").Select()
